$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.751.76"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "1.926.66"
$ws.Range("E3").Value = "  -1.48%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "241.99"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "0.4848"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").Value = "0.2929"
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("D9").Value = "0.06807"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("E11").Value = "  -1.10%  "
$ws.Range("D12").Value = "1.937.74"
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("D13").Value = "0.07764"
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("D14").Value = "5.311"
$ws.Range("E14").Value = "  -2.67%  "
$ws.Range("D15").Value = "0.6951"
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("D16").Value = "274.50"
$ws.Range("E16").Value = "  -3.21%  "
$ws.Range("D17").Value = "30.753.79"
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("D18").Value = "0.000007647"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("D20").Value = "12.92"
$ws.Range("E20").Value = "  -1.97%  "
$ws.Range("D21").Value = "5.548"
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").Value = "6.446"
$ws.Range("E23").Value = "  -0.72%  "
$ws.Range("D24").Value = "9.833"
$ws.Range("E24").Value = "  +0.28%  "
$ws.Range("D25").Value = "164.43"
$ws.Range("E25").Value = "  -3.25%  "
$ws.Range("D26").Value = "19.40"
$ws.Range("E26").Value = "  -2.76%  "
$ws.Range("D27").Value = "2.143"
$ws.Range("E27").Value = "  -2.53%  "
$ws.Range("D28").Value = "0.1036"
$ws.Range("E28").Value = "  -1.69%  "
$ws.Range("E29").Value = "  -1.75%  "
$ws.Range("D30").Value = "4.558"
$ws.Range("E30").Value = "  -1.13%  "
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("D32").Value = "4.339"
$ws.Range("E32").Value = "  -2.26%  "
$ws.Range("E33").Value = "  -1.14%  "
$ws.Range("D34").Value = "0.7541"
$ws.Range("E34").Value = "  -1.16%  "
$ws.Range("E35").Value = "  -2.76%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").Value = "2.718"
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("D38").Value = "0.01983"
$ws.Range("E38").Value = "  -1.07%  "
$ws.Range("D39").Value = "2.647"
$ws.Range("E39").Value = "  -2.12%  "
$ws.Range("D40").Value = "6.437"
$ws.Range("E40").Value = "  -1.54%  "
$ws.Range("D41").Value = "77.15"
$ws.Range("E41").Value = "  +2.85%  "
$ws.Range("D42").Value = "2.051"
$ws.Range("E42").Value = "  -2.30%  "
$ws.Range("D43").Value = "0.8806"
$ws.Range("E43").Value = "  -0.95%  "
$ws.Range("D44").Value = "0.4413"
$ws.Range("E44").Value = "  -1.02%  "
$ws.Range("D45").Value = "107.27"
$ws.Range("E45").Value = "  -1.80%  "
$ws.Range("D46").Value = "7.826"
$ws.Range("E46").Value = "  -4.35%  "
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").Value = "980.04"
$ws.Range("E48").Value = "  -2.18%  "
$ws.Range("D49").Value = "36.01"
$ws.Range("E49").Value = "  +0.75%  "
$ws.Range("D50").Value = "0.1232"
$ws.Range("E50").Value = "  -1.86%  "
$ws.Range("D51").Value = "9.114"
$ws.Range("E51").Value = "  -2.40%  "
